# Weekly update: insert a new daily price record for Zanahoria
# (Agricola del Norte S.A. de Arica) as row 105, pushing the existing
# rows 105-216 down to 106-217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 105 (shifts rows 105..216 -> 106..217)
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A105").Value = 1
$ws.Range("B105").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C105").Value = "Arica y Parinacota"
$ws.Range("D105").Value = 44539
$ws.Range("E105").Value = 15
$ws.Range("F105").Value = 100114013
$ws.Range("G105").Value = "Zanahoria"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 90
$ws.Range("K105").Value = 18000
$ws.Range("L105").Value = 20000
$ws.Range("M105").Value = 19000
$ws.Range("N105").Value = "$/saco 25 kilos"
$ws.Range("O105").Value = "Provincia de Calama"
$ws.Range("P105").Value = 760
$ws.Range("Q105").Value = 25
$ws.Range("R105").Value = "Hortaliza"
